$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$nl = [char]11

$cell = $t.Cell(1, 1)
$cell.Range.Text = '74 x 35' + $nl + '  3    5' + $nl + '  ----' + $nl + '7|    |' + $nl + '4|    |'

$cell = $t.Cell(1, 2)
$cell.Range.Text = '50 x 42' + $nl + '  4    2' + $nl + '  ----' + $nl + '5|    |' + $nl + '0|    |'

$cell = $t.Cell(1, 3)
$cell.Range.Text = '66 x 56' + $nl + '  5    6' + $nl + '  ----' + $nl + '6|    |' + $nl + '6|    |'

$cell = $t.Cell(2, 1)
$cell.Range.Text = '36 x 70' + $nl + '  7    0' + $nl + '  ----' + $nl + '3|    |' + $nl + '6|    |'

$cell = $t.Cell(2, 2)
$cell.Range.Text = '69 x 29' + $nl + '  2    9' + $nl + '  ----' + $nl + '6|    |' + $nl + '9|    |'

$cell = $t.Cell(2, 3)
$cell.Range.Text = '15 x 71' + $nl + '  7    1' + $nl + '  ----' + $nl + '1|    |' + $nl + '5|    |'

$cell = $t.Cell(3, 1)
$cell.Range.Text = '48 x 66' + $nl + '  6    6' + $nl + '  ----' + $nl + '4|    |' + $nl + '8|    |'

$cell = $t.Cell(3, 2)
$cell.Range.Text = '38 x 32' + $nl + '  3    2' + $nl + '  ----' + $nl + '3|    |' + $nl + '8|    |'

$cell = $t.Cell(3, 3)
$cell.Range.Text = '59 x 28' + $nl + '  2    8' + $nl + '  ----' + $nl + '5|    |' + $nl + '9|    |'

$cell = $t.Cell(4, 1)
$cell.Range.Text = '37 x 63' + $nl + '  6    3' + $nl + '  ----' + $nl + '3|    |' + $nl + '7|    |'

$cell = $t.Cell(4, 2)
$cell.Range.Text = '33 x 22' + $nl + '  2    2' + $nl + '  ----' + $nl + '3|    |' + $nl + '3|    |'

$cell = $t.Cell(4, 3)
$cell.Range.Text = '28 x 59' + $nl + '  5    9' + $nl + '  ----' + $nl + '2|    |' + $nl + '8|    |'

$cell = $t.Cell(5, 1)
$cell.Range.Text = '18 x 28' + $nl + '  2    8' + $nl + '  ----' + $nl + '1|    |' + $nl + '8|    |'

$cell = $t.Cell(5, 2)
$cell.Range.Text = '56 x 87' + $nl + '  8    7' + $nl + '  ----' + $nl + '5|    |' + $nl + '6|    |'

$cell = $t.Cell(5, 3)
$cell.Range.Text = '50 x 92' + $nl + '  9    2' + $nl + '  ----' + $nl + '5|    |' + $nl + '0|    |'
